$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; existing rows 19..45 shift down to 20..46,
# and the sheet dimension grows from R45 to R46 automatically.
$ws.Rows(19).Insert()

# Populate the newly inserted row 19 with this week's new record.
$ws.Cells(19,1).Value  = 1
$ws.Cells(19,2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells(19,3).Value  = 'Arica y Parinacota'
$ws.Cells(19,4).Value  = 45096
$ws.Cells(19,5).Value  = 15
$ws.Cells(19,6).Value  = 100112044
$ws.Cells(19,7).Value  = 'Perejil'
$ws.Cells(19,8).Value  = 'Sin especificar'
$ws.Cells(19,9).Value  = 'Primera'
$ws.Cells(19,10).Value = 250
$ws.Cells(19,11).Value = 900
$ws.Cells(19,12).Value = 1000
$ws.Cells(19,13).Value = 950
$ws.Cells(19,14).Value = '$/atado 1,5 a 2 kilos'
$ws.Cells(19,15).Value = 'Región de Arica y Parinacota'
$ws.Cells(19,16).Value = 475
$ws.Cells(19,17).Value = 2
$ws.Cells(19,18).Value = 'Hortaliza'
